$d = $word.ActiveDocument

# The page removed its trailing "Ver no Jupiter / Salvar..." and
# "(c) 2020 ... Jekyll ..." footer boilerplate, along with the blank
# paragraph that separated it from the preceding "Requisitos" text
# ("LOQ4038: Quimica Organica II (Requisito fraco)"). Locate the
# "Ver no Jupiter" paragraph, then widen the range one paragraph in
# each direction so it spans [blank paragraph] [Ver no Jupiter ...]
# [(c) 2020 ...] and delete that whole range in one shot, leaving the
# "LOQ4038..." paragraph directly followed by the document's final
# blank paragraph / page-break paragraph.

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ver no Jupiter*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $startPara = $d.Paragraphs.Item($targetIndex - 1)
    $endPara = $d.Paragraphs.Item($targetIndex + 1)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
